$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("payment-request")
$ws.Name = "testSheet"

$ws.Range("A1").Value = "insert"
$ws.Range("B1").Value = "dobpaymentrequest"
$ws.Range("A6").Value = "insert"
$ws.Range("B6").Value = "dobpaymentrequest"
$ws.Range("A11").Value = "insert"
$ws.Range("B11").Value = "dobpaymentrequest"
$ws.Range("K8").Value = "null"
$ws.Range("J13").Value = "'" + [char]8217
